$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the header row (row 1), before the
# existing first data row (001761119 / BLUEMETRIX).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Row 2: 000806386 / FERNANDA / 1752706.81
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "000806386"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = "FERNANDA"
$ws.Cells.Item(2, 3).Value = 1752706.81

# Row 3: 004228090 / GUSTAVO / 700185.23
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "004228090"
$ws.Cells.Item(3, 1).ClearFormats()
$ws.Cells.Item(3, 2).Value = "GUSTAVO"
$ws.Cells.Item(3, 3).Value = 700185.23

# Remove the old last data row (004361159 / HFR / -149433.71), which is
# now on row 241 after the two inserts above (was row 239).
$ws.Rows.Item(241).Delete()
